$wb = $excel.ActiveWorkbook

$wsCompletion = $wb.Worksheets.Item("Completion Dates")
$wsResults = $wb.Worksheets.Item("Results")
$wsIntervention = $wb.Worksheets.Item("intervention")

# Delete the duplicate row (row 740) on both sheets; this shifts all rows below up by one.
$wsCompletion.Rows.Item(740).Delete()
$wsResults.Rows.Item(740).Delete()

# Restore view/selection state seen in the target workbook.
$wsIntervention.Range("A889").Select()

$wsCompletion.Range("A739").Select()
$wsResults.Range("A739").Select()
